# Adding more countries to compare, translating, updating data
# - Append a new date row (2020-06-20) to "Confirmados" and "Mortes"
# - Correct the last-day-of-range figures that were carried over as
#   provisional values in row 116 for both sheets

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Confirmados")
$ws2 = $wb.Worksheets.Item("Mortes")

# ---------------------------------------------------------------------
# Confirmados (sheet1): corrections to the existing last row (116)
# ---------------------------------------------------------------------
$ws1.Range("K116").Value = 31760
$ws1.Range("N116").Value = 19206
$ws1.Range("Q116").Value = 4990
$ws1.Range("V116").Value = 8203
$ws1.Range("Z116").Value = 8725
$ws1.Range("AB116").Value = 8037

# ---------------------------------------------------------------------
# Confirmados (sheet1): new row for 2020-06-20
# ---------------------------------------------------------------------
$row1 = @{
    "A" = "'2020-06-20"
    "B" = 215793
    "C" = 34238
    "D" = 95537
    "E" = 45304
    "F" = 32235
    "G" = 27698
    "H" = 27305
    "I" = 19138
    "J" = 15418
    "K" = 36521
    "L" = 51118
    "M" = 13821
    "N" = 19206
    "O" = 17108
    "P" = 62902
    "Q" = 4990
    "R" = 18449
    "S" = 92397
    "T" = 11263
    "U" = 84654
    "V" = 8203
    "W" = 14264
    "X" = 21574
    "Y" = 69673
    "Z" = 9262
    "AA" = 14952
    "AB" = 8037
}
foreach ($col in $row1.Keys) {
    $ws1.Range("$col" + "117").Value = $row1[$col]
}

# ---------------------------------------------------------------------
# Mortes (sheet2): corrections to existing rows
# ---------------------------------------------------------------------
$ws2.Range("AB114").Value = 224

$ws2.Range("F116").Value = 396
$ws2.Range("K116").Value = 709
$ws2.Range("N116").Value = 697
$ws2.Range("Q116").Value = 43
$ws2.Range("V116").Value = 167
$ws2.Range("Z116").Value = 322
$ws2.Range("AB116").Value = 245

# ---------------------------------------------------------------------
# Mortes (sheet2): new row for 2020-06-20
# ---------------------------------------------------------------------
$row2 = @{
    "A" = "'2020-06-20"
    "B" = 12494
    "C" = 1297
    "D" = 8824
    "E" = 1350
    "F" = 405
    "G" = 866
    "H" = 636
    "I" = 430
    "J" = 295
    "K" = 743
    "L" = 4148
    "M" = 437
    "N" = 697
    "O" = 237
    "P" = 2650
    "Q" = 43
    "R" = 448
    "S" = 5520
    "T" = 297
    "U" = 4583
    "V" = 167
    "W" = 485
    "X" = 353
    "Y" = 1684
    "Z" = 341
    "AA" = 409
    "AB" = 245
}
foreach ($col in $row2.Keys) {
    $ws2.Range("$col" + "117").Value = $row2[$col]
}

# ---------------------------------------------------------------------
# Tidy up: the leading apostrophe used to force the date-like label to be
# stored as text (matching columns A1:A116) leaves a transient
# quote-prefixed style on A117; reset both new date cells back to the
# plain "Normal" style so they match the rest of column A.
# ---------------------------------------------------------------------
$ws1.Range("A117").Style = "Normal"
$ws2.Range("A117").Style = "Normal"
